$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A (TabName) and a new row 3 (FilesTab) ---
$ws.Columns("A").Insert()
$ws.Rows("3").Insert()

$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.pubmed_id IN ['31504139'] 
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.pubmed_id IN ['31504139']
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@
$statsQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
  WHERE a.pubmed_id IN ['31504139']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# --- Column A (TabName, CasesTab, FilesTab) filled first so the new shared
#     strings are appended in the same order as the authored workbook ---
$ws.Range("A1").Value2 = "TabName"
$ws.Range("A2").Value2 = "CasesTab"
$ws.Range("A3").Value2 = "FilesTab"

# --- Query columns (B = query, C = StatQuery) ---
$ws.Range("B2").Value2 = $casesQuery
$ws.Range("B3").Value2 = $filesQuery
$ws.Range("C2").Value2 = $statsQuery
$ws.Range("C3").Value2 = $statsQuery

# --- Remaining row 1 headers ---
$ws.Range("B1").Value2 = "query"
$ws.Range("C1").Value2 = "StatQuery"
$ws.Range("D1").Value2 = "dbExcel"
$ws.Range("E1").Value2 = "WebExcel"

# --- Remaining data columns (D/E) ---
$ws.Range("D2").Value2 = "TC01_Trials_Filter_PubmedID-315_Neo4jData.xlsx"
$ws.Range("E2").Value2 = "TC01_Trials_Filter_PubmedID-315_WebData.xlsx"
$ws.Range("D3").Value2 = "TC01_Trials_Filter_PubmedID-315_Neo4jData.xlsx"
$ws.Range("E3").Value2 = "TC01_Trials_Filter_PubmedID-315_WebData.xlsx"

# --- Styles: wrap text on B2, C2, B3, C3 ---
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Row heights ---
$ws.Rows("2").RowHeight = 195
$ws.Rows("3").RowHeight = 409.5

# --- Column widths (best achievable precision in this runtime) ---
$ws.Columns("A").ColumnWidth = 8
$ws.Columns("B").ColumnWidth = 75
$ws.Columns("C").ColumnWidth = 75
$ws.Columns("D").ColumnWidth = 69.5
$ws.Columns("E").ColumnWidth = 27.666666666666668

# --- Selection / view ---
$ws.Range("C3").Select() | Out-Null
